$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 123.6
$ws.Range("I4").Value = 123.6
$ws.Range("K4").Value = 123.6
$ws.Range("M4").Value = -9.599999999999994
$ws.Range("H43").Value = 1135
$ws.Range("I43").Value = 1095
$ws.Range("K43").Value = 1095
$ws.Range("M43").Value = -1026
$ws.Range("H55").Value = 100.92308
$ws.Range("I55").Value = 75.7
$ws.Range("J55").Value = 185
$ws.Range("K55").Value = 75.7
$ws.Range("L55").Value = 185
$ws.Range("M55").Value = 138.3
$ws.Range("N55").Value = -613
$ws.Range("H76").Value = 6794.125
$ws.Range("I76").Value = 5911.778
$ws.Range("K76").Value = 5911.778
$ws.Range("M76").Value = -5596.778
$ws.Range("H79").Value = 6794.125
$ws.Range("I79").Value = 5911.778
$ws.Range("K79").Value = 5911.778
$ws.Range("M79").Value = -4819.778
$ws.Range("H132").Value = 3492.7144
$ws.Range("I132").Value = 3530.6155
$ws.Range("K132").Value = 10591.8465
$ws.Range("M132").Value = -8061.8465

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 33348504
$ws.Range("I132").Value = 6349.654
$ws.Range("K132").Value = 19048.962
$ws.Range("M132").Value = -16518.962

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2619.9285
$ws.Range("I31").Value = 1713.05
$ws.Range("J31").Value = 4887.125
$ws.Range("K31").Value = 1713.05
$ws.Range("L31").Value = 4887.125
$ws.Range("M31").Value = -1418.05
$ws.Range("N31").Value = -5477.125
$ws.Range("H34").Value = 2619.9285
$ws.Range("I34").Value = 1713.05
$ws.Range("J34").Value = 4887.125
$ws.Range("K34").Value = 1713.05
$ws.Range("L34").Value = 4887.125
$ws.Range("M34").Value = -1511.05
$ws.Range("N34").Value = -5291.125
$ws.Range("H48").Value = 36000
$ws.Range("J48").Value = 36000
$ws.Range("L48").Value = 36000
$ws.Range("N48").Value = -36952
$ws.Range("H86").Value = 6326.5
$ws.Range("I86").Value = 5900.1
$ws.Range("K86").Value = 5900.1
$ws.Range("M86").Value = -4777.1
$ws.Range("H89").Value = 6326.5
$ws.Range("I89").Value = 5900.1
$ws.Range("K89").Value = 29500.5
$ws.Range("M89").Value = -23884.5
$ws.Range("H105").Value = 8940.823
$ws.Range("I105").Value = 1908.5454
$ws.Range("K105").Value = 1908.5454
$ws.Range("M105").Value = -161.5454
$ws.Range("H132").Value = 2660.9614
$ws.Range("J132").Value = 5071
$ws.Range("L132").Value = 15213
$ws.Range("N132").Value = -20273

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 1365.1666
$ws.Range("I22").Value = 1288.7
$ws.Range("J22").Value = 1747.5
$ws.Range("K22").Value = 3866.1
$ws.Range("L22").Value = 5242.5
$ws.Range("M22").Value = -3697.1
$ws.Range("N22").Value = -5580.5
$ws.Range("H27").Value = 1365.1666
$ws.Range("I27").Value = 1288.7
$ws.Range("J27").Value = 1747.5
$ws.Range("K27").Value = 3866.1
$ws.Range("L27").Value = 5242.5
$ws.Range("M27").Value = -3764.1
$ws.Range("N27").Value = -5446.5
$ws.Range("H34").Value = 539
$ws.Range("J34").Value = 1498
$ws.Range("L34").Value = 4494
$ws.Range("N34").Value = -4662
$ws.Range("H39").Value = 2493
$ws.Range("H55").Value = 1003033.1
$ws.Range("J55").Value = 6000
$ws.Range("L55").Value = 18000
$ws.Range("N55").Value = -18354
$ws.Range("H139").Value = 3436.8572
$ws.Range("I139").Value = 3176.5
$ws.Range("K139").Value = 9529.5
$ws.Range("M139").Value = -4389.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 25366666
$ws.Range("J7").Value = 50000
$ws.Range("L7").Value = 50000
$ws.Range("N7").Value = -50224
$ws.Range("H8").Value = 25366666
$ws.Range("J8").Value = 50000
$ws.Range("L8").Value = 50000
$ws.Range("N8").Value = -50278
$ws.Range("H11").Value = 1020067.75
$ws.Range("I11").Value = 2220199
$ws.Range("J11").Value = 19958.334
$ws.Range("K11").Value = 2220199
$ws.Range("L11").Value = 19958.334
$ws.Range("M11").Value = -2220060
$ws.Range("N11").Value = -20236.334
$ws.Range("H123").Value = 51262.5
$ws.Range("I123").Value = 57600
$ws.Range("K123").Value = 57600
$ws.Range("M123").Value = -55150
$ws.Range("H126").Value = 4598.143
$ws.Range("I126").Value = 4297
$ws.Range("J126").Value = 4999.6665
$ws.Range("K126").Value = 12891
$ws.Range("L126").Value = 14998.9995
$ws.Range("M126").Value = -10421
$ws.Range("N126").Value = -19938.9995
$ws.Range("H132").Value = 4301.4287
$ws.Range("I132").Value = 3698.9443
$ws.Range("J132").Value = 7916.3335
$ws.Range("K132").Value = 11096.8329
$ws.Range("L132").Value = 23749.0005
$ws.Range("M132").Value = -8566.832900000001
$ws.Range("N132").Value = -28809.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 2750
$ws.Range("I3").Value = 4500
$ws.Range("J3").Value = 1000
$ws.Range("K3").Value = 4500
$ws.Range("L3").Value = 1000
$ws.Range("M3").Value = -4388
$ws.Range("N3").Value = -1224
$ws.Range("H14").Value = 25500
$ws.Range("J14").Value = 25500
$ws.Range("L14").Value = 25500
$ws.Range("N14").Value = -25844
$ws.Range("H15").Value = 2750
$ws.Range("I15").Value = 4500
$ws.Range("J15").Value = 1000
$ws.Range("K15").Value = 4500
$ws.Range("L15").Value = 1000
$ws.Range("M15").Value = -4330
$ws.Range("N15").Value = -1340
$ws.Range("H55").Value = 134.6
$ws.Range("I55").Value = 134.6
$ws.Range("K55").Value = 134.6
$ws.Range("M55").Value = 38.40000000000001
$ws.Range("H68").Value = 4959.5
$ws.Range("I68").Value = 5112.6665
$ws.Range("K68").Value = 5112.6665
$ws.Range("M68").Value = -4363.6665
$ws.Range("H71").Value = 4959.5
$ws.Range("I71").Value = 5112.6665
$ws.Range("K71").Value = 25563.3325
$ws.Range("M71").Value = -21819.3325
$ws.Range("H136").Value = 1820626.9
$ws.Range("I136").Value = 4001919.8
$ws.Range("K136").Value = 12005759.4
$ws.Range("M136").Value = -12003209.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 11000
$ws.Range("J18").Value = 11000
$ws.Range("L18").Value = 11000
$ws.Range("N18").Value = -11346
$ws.Range("H55").Value = 30029
$ws.Range("J55").Value = 30029
$ws.Range("L55").Value = 30029
$ws.Range("N55").Value = -30583
